# Update symbol-list values (price/volume columns are stored as literal
# text in the source sheet, so numeric-looking values are entered with a
# leading apostrophe to keep Excel from auto-converting them to numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.43"
$ws.Range("E2").Value = "'-1.17%"
$ws.Range("D3").Value = "'43.87"
$ws.Range("E3").Value = "'5.74%"
$ws.Range("D4").Value = "'5.486"
$ws.Range("E4").Value = "'-2.31%"
$ws.Range("D5").Value = "'0.08152"
$ws.Range("E5").Value = "'-2.25%"
$ws.Range("D6").Value = "'8.705"
$ws.Range("E6").Value = "'-0.89%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.932"
$ws.Range("E7").Value = "'-2.84%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.325"
$ws.Range("E8").Value = "'-3.48%"
$ws.Range("E9").Value = "'-3.04%"
$ws.Range("D10").Value = "'0.9412"
$ws.Range("E10").Value = "'1.74%"
$ws.Range("D11").Value = "'0.1182"
$ws.Range("E11").Value = "'-8.38%"
$ws.Range("D12").Value = "'0.1898"
$ws.Range("E12").Value = "'-3.73%"
$ws.Range("D13").Value = "'0.09845"
$ws.Range("E13").Value = "'3.03%"
$ws.Range("E14").Value = "'7.24%"
$ws.Range("E15").Value = "'0.91%"
$ws.Range("D16").Value = "'0.001313"
$ws.Range("E16").Value = "'1.05%"
$ws.Range("D17").Value = "'0.006112"
$ws.Range("E17").Value = "'1.73%"
$ws.Range("D18").Value = "'3.553"
$ws.Range("E18").Value = "'3.27%"
$ws.Range("E19").Value = "'-0.52%"
$ws.Range("D20").Value = "'8.731"
$ws.Range("E20").Value = "'2.39%"
$ws.Range("D21").Value = "'0.1351"
$ws.Range("E21").Value = "'-1.42%"
$ws.Range("D22").Value = "'0.2500"
$ws.Range("E22").Value = "'2.32%"
$ws.Range("D23").Value = "'0.04386"
$ws.Range("E23").Value = "'-0.30%"
$ws.Range("E24").Value = "'-2.53%"
$ws.Range("D25").Value = "'0.004338"
$ws.Range("E25").Value = "'-1.12%"
$ws.Range("D26").Value = "'0.0001237"
$ws.Range("E26").Value = "'3.13%"
$ws.Range("D27").Value = "'0.0004012"
$ws.Range("E27").Value = "'31.73%"
$ws.Range("D39").Value = "'0.02667"
$ws.Range("E39").Value = "'-5.12%"
$ws.Range("D40").Value = "'0.05630"
$ws.Range("E40").Value = "'2.16%"
$ws.Range("D41").Value = "'0.007882"
$ws.Range("E41").Value = "'-0.93%"
$ws.Range("D42").Value = "'0.009789"
$ws.Range("E42").Value = "'5.14%"
$ws.Range("D43").Value = "'0.1409"
$ws.Range("E43").Value = "'-1.97%"
$ws.Range("D44").Value = "'0.002108"
$ws.Range("E44").Value = "'-1.43%"
$ws.Range("D45").Value = "'0.009629"
$ws.Range("E45").Value = "'-12.87%"
$ws.Range("D46").Value = "'0.00007060"
$ws.Range("E46").Value = "'0.01%"
$ws.Range("D47").Value = "'0.00000000754"
$ws.Range("E47").Value = "'0.62%"
$ws.Range("D48").Value = "'0.003456"
$ws.Range("E48").Value = "'6.78%"
$ws.Range("D49").Value = "'0.002283"
$ws.Range("E49").Value = "'0.28%"
$ws.Range("D50").Value = "'0.00002112"
$ws.Range("E50").Value = "'0.62%"
$ws.Range("D51").Value = "'0.0002011"
$ws.Range("E51").Value = "'0.62%"

Write-Host "Applied all changes"